$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D51:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.585.88"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.717.27"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D5").Value = "240.58"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.4921"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").Value = "0.2596"
$ws.Range("D9").Value = "0.06204"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "1.729.44"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "0.06997"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "15.71"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "0.6073"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "4.478"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "76.66"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "0.9994"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "26.443.98"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "0.9988"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "0.000007133"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").Value = "1.952.73"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "4.408"
$ws.Range("E22").Value = "  -3.28%  "
$ws.Range("D23").Value = "8.491"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").Value = "5.085"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").Value = "137.47"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "15.24"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("D29").Value = "105.59"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "3.913"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").Value = "0.07937"
$ws.Range("D32").Value = "3.638"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").Value = "0.04509"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.615"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "0.9974"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6251"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "0.9383"
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "2.008"
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.406"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "0.9986"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01500"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "99.55"
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.508"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.3834"
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "6.937"
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1154"
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05371"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "7.749"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "30.09"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "51.43"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.221"
$ws.Range("E51").Value = "  -2.47%  "
